$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.4607133333333334
$ws.Cells.Item(2, 8).Value = 1.38214
$ws.Cells.Item(2, 9).Value = 0.1068584005969239
$ws.Cells.Item(2, 10).Value = 0.1068584005969239
$ws.Cells.Item(2, 15).Value = 0.02266023449704293
$ws.Cells.Item(2, 16).Value = 0.02266023449704293
$ws.Cells.Item(2, 17).Value = 0.07254192504222222
$ws.Cells.Item(2, 18).Value = 0.6528773253800001
$ws.Cells.Item(2, 19).Value = 0.002421436415505247
$ws.Cells.Item(2, 20).Value = 0.002421436415505247
$ws.Cells.Item(3, 7).Value = 0.4607133333333334
$ws.Cells.Item(3, 8).Value = 1.38214
$ws.Cells.Item(3, 9).Value = 0.1068584005969239
$ws.Cells.Item(3, 10).Value = 0.1068584005969239
$ws.Cells.Item(3, 13).Value = 6.739756333333333
$ws.Cells.Item(3, 15).Value = 0.9699521281096917
$ws.Cells.Item(3, 16).Value = 0.9699521281096917
$ws.Cells.Item(3, 17).Value = 3.105095606184445
$ws.Cells.Item(3, 18).Value = 27.94586045566
$ws.Cells.Item(3, 19).Value = 0.1036475330653843
$ws.Cells.Item(3, 20).Value = 0.1036475330653843
$ws.Cells.Item(4, 7).Value = 0.4607133333333334
$ws.Cells.Item(4, 8).Value = 1.38214
$ws.Cells.Item(4, 9).Value = 0.1068584005969239
$ws.Cells.Item(4, 10).Value = 0.1068584005969239
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.05133333333333333
$ws.Cells.Item(4, 14).Value = 0.154
$ws.Cells.Item(4, 15).Value = 0.00738763739326543
$ws.Cells.Item(4, 16).Value = 0.00738763739326543
$ws.Cells.Item(4, 17).Value = 0.02364995111111111
$ws.Cells.Item(4, 18).Value = 0.21284956
$ws.Cells.Item(4, 19).Value = 0.0007894311160343717
$ws.Cells.Item(4, 20).Value = 0.0007894311160343717
$ws.Cells.Item(5, 9).Value = 0.2623141252264423
$ws.Cells.Item(5, 10).Value = 0.2623141252264423
$ws.Cells.Item(5, 15).Value = 0.02266023449704293
$ws.Cells.Item(5, 16).Value = 0.02266023449704293
$ws.Cells.Item(5, 19).Value = 0.005944099589517866
$ws.Cells.Item(5, 20).Value = 0.005944099589517866
$ws.Cells.Item(6, 9).Value = 0.2623141252264423
$ws.Cells.Item(6, 10).Value = 0.2623141252264423
$ws.Cells.Item(6, 13).Value = 6.739756333333333
$ws.Cells.Item(6, 15).Value = 0.9699521281096917
$ws.Cells.Item(6, 16).Value = 0.9699521281096917
$ws.Cells.Item(6, 17).Value = 7.622334164939667
$ws.Cells.Item(6, 18).Value = 68.601007484457
$ws.Cells.Item(6, 19).Value = 0.2544321439966198
$ws.Cells.Item(6, 20).Value = 0.2544321439966198
$ws.Cells.Item(7, 9).Value = 0.2623141252264423
$ws.Cells.Item(7, 10).Value = 0.2623141252264423
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.05133333333333333
$ws.Cells.Item(7, 14).Value = 0.154
$ws.Cells.Item(7, 15).Value = 0.00738763739326543
$ws.Cells.Item(7, 16).Value = 0.00738763739326543
$ws.Cells.Item(7, 17).Value = 0.05805548466666667
$ws.Cells.Item(7, 18).Value = 0.522499362
$ws.Cells.Item(7, 19).Value = 0.001937881640304576
$ws.Cells.Item(7, 20).Value = 0.001937881640304575
$ws.Cells.Item(8, 7).Value = 1.908186666666666
$ws.Cells.Item(8, 8).Value = 5.724559999999999
$ws.Cells.Item(8, 9).Value = 0.4425870937250397
$ws.Cells.Item(8, 10).Value = 0.4425870937250397
$ws.Cells.Item(8, 15).Value = 0.02266023449704293
$ws.Cells.Item(8, 16).Value = 0.02266023449704293
$ws.Cells.Item(8, 17).Value = 0.3004548037244444
$ws.Cells.Item(8, 18).Value = 2.70409323352
$ws.Cells.Item(8, 19).Value = 0.01002912732917412
$ws.Cells.Item(8, 20).Value = 0.01002912732917412
$ws.Cells.Item(9, 7).Value = 1.908186666666666
$ws.Cells.Item(9, 8).Value = 5.724559999999999
$ws.Cells.Item(9, 9).Value = 0.4425870937250397
$ws.Cells.Item(9, 10).Value = 0.4425870937250397
$ws.Cells.Item(9, 13).Value = 6.739756333333333
$ws.Cells.Item(9, 15).Value = 0.9699521281096917
$ws.Cells.Item(9, 16).Value = 0.9699521281096917
$ws.Cells.Item(9, 17).Value = 12.86071317184889
$ws.Cells.Item(9, 18).Value = 115.74641854664
$ws.Cells.Item(9, 19).Value = 0.4292882934324859
$ws.Cells.Item(9, 20).Value = 0.4292882934324859
$ws.Cells.Item(10, 7).Value = 1.908186666666666
$ws.Cells.Item(10, 8).Value = 5.724559999999999
$ws.Cells.Item(10, 9).Value = 0.4425870937250397
$ws.Cells.Item(10, 10).Value = 0.4425870937250397
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.05133333333333333
$ws.Cells.Item(10, 14).Value = 0.154
$ws.Cells.Item(10, 15).Value = 0.00738763739326543
$ws.Cells.Item(10, 16).Value = 0.00738763739326543
$ws.Cells.Item(10, 17).Value = 0.09795358222222221
$ws.Cells.Item(10, 18).Value = 0.8815822399999999
$ws.Cells.Item(10, 19).Value = 0.003269672963379775
$ws.Cells.Item(10, 20).Value = 0.003269672963379775
$ws.Cells.Item(11, 7).Value = 0.632459
$ws.Cells.Item(11, 8).Value = 1.897377
$ws.Cells.Item(11, 9).Value = 0.1466932955774304
$ws.Cells.Item(11, 10).Value = 0.1466932955774304
$ws.Cells.Item(11, 15).Value = 0.02266023449704293
$ws.Cells.Item(11, 16).Value = 0.02266023449704293
$ws.Cells.Item(11, 17).Value = 0.09958425348433333
$ws.Cells.Item(11, 18).Value = 0.896258281359
$ws.Cells.Item(11, 19).Value = 0.003324104476928603
$ws.Cells.Item(11, 20).Value = 0.003324104476928603
$ws.Cells.Item(12, 7).Value = 0.632459
$ws.Cells.Item(12, 8).Value = 1.897377
$ws.Cells.Item(12, 9).Value = 0.1466932955774304
$ws.Cells.Item(12, 10).Value = 0.1466932955774304
$ws.Cells.Item(12, 13).Value = 6.739756333333333
$ws.Cells.Item(12, 15).Value = 0.9699521281096917
$ws.Cells.Item(12, 16).Value = 0.9699521281096917
$ws.Cells.Item(12, 17).Value = 4.262619550823667
$ws.Cells.Item(12, 18).Value = 38.363575957413
$ws.Cells.Item(12, 19).Value = 0.1422854742247526
$ws.Cells.Item(12, 20).Value = 0.1422854742247526
$ws.Cells.Item(13, 7).Value = 0.632459
$ws.Cells.Item(13, 8).Value = 1.897377
$ws.Cells.Item(13, 9).Value = 0.1466932955774304
$ws.Cells.Item(13, 10).Value = 0.1466932955774304
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.05133333333333333
$ws.Cells.Item(13, 14).Value = 0.154
$ws.Cells.Item(13, 15).Value = 0.00738763739326543
$ws.Cells.Item(13, 16).Value = 0.00738763739326543
$ws.Cells.Item(13, 17).Value = 0.03246622866666667
$ws.Cells.Item(13, 18).Value = 0.292196058
$ws.Cells.Item(13, 19).Value = 0.001083716875749163
$ws.Cells.Item(13, 20).Value = 0.001083716875749163
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.1791276666666667
$ws.Cells.Item(14, 8).Value = 0.5373830000000001
$ws.Cells.Item(14, 9).Value = 0.04154708487416379
$ws.Cells.Item(14, 10).Value = 0.04154708487416379
$ws.Cells.Item(14, 15).Value = 0.02266023449704293
$ws.Cells.Item(14, 16).Value = 0.02266023449704293
$ws.Cells.Item(14, 17).Value = 0.02820466617344445
$ws.Cells.Item(14, 18).Value = 0.253841995561
$ws.Cells.Item(14, 19).Value = 0.0009414666859170968
$ws.Cells.Item(14, 20).Value = 0.0009414666859170967
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.1791276666666667
$ws.Cells.Item(15, 8).Value = 0.5373830000000001
$ws.Cells.Item(15, 9).Value = 0.04154708487416379
$ws.Cells.Item(15, 10).Value = 0.04154708487416379
$ws.Cells.Item(15, 13).Value = 6.739756333333333
$ws.Cells.Item(15, 15).Value = 0.9699521281096917
$ws.Cells.Item(15, 16).Value = 0.9699521281096917
$ws.Cells.Item(15, 17).Value = 1.207276825891889
$ws.Cells.Item(15, 18).Value = 10.865491433027
$ws.Cells.Item(15, 19).Value = 0.04029868339044915
$ws.Cells.Item(15, 20).Value = 0.04029868339044915
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.1791276666666667
$ws.Cells.Item(16, 8).Value = 0.5373830000000001
$ws.Cells.Item(16, 9).Value = 0.04154708487416379
$ws.Cells.Item(16, 10).Value = 0.04154708487416379
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.05133333333333333
$ws.Cells.Item(16, 14).Value = 0.154
$ws.Cells.Item(16, 15).Value = 0.00738763739326543
$ws.Cells.Item(16, 16).Value = 0.00738763739326543
$ws.Cells.Item(16, 17).Value = 0.009195220222222224
$ws.Cells.Item(16, 18).Value = 0.08275698200000001
$ws.Cells.Item(16, 19).Value = 0.000306934797797545
$ws.Cells.Item(16, 20).Value = 0.0003069347977975449
